$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update handoff/handback datetimes for row 3
# (6603f24f-90f7-44ae-8799-dd3d11673c6d.md entry)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-23 04:02:52"
$wsZh.Range("H3").Value = "2016-03-23 04:03:36"

# "de-de" sheet: update handoff/handback datetimes for row 3
# (6603f24f-90f7-44ae-8799-dd3d11673c6d.md entry)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-23 04:03:01"
$wsDe.Range("H3").Value = "2016-03-23 04:03:51"
